# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") held values like "4-27-2012-13" which need to
# become the corrected date-as-text "2013-04-27".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"

$range = $ws.Range($col + $firstRow + ":" + $col + $lastRow)
# Force text storage so Excel doesn't reinterpret the new value as a date.
$range.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range($col + $row).Value = "2013-04-27"
}
